$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.100.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "'3.055.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D5").Value = "'390.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "'100.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("D7").Value = "'0.532"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").Value = "'3.536.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "'18.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "'7.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "'3.031.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "'10.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "'51.094.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").Value = "'12.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "'69.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "'263.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("D26").Value = "'7.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.97%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "'7.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.99%  "
$ws.Range("E30").Value = "  -5.77%  "
$ws.Range("E31").Value = "  -3.17%  "
$ws.Range("D32").Value = "'10.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("D33").Value = "'0.0487"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.60%  "
$ws.Range("D34").Value = "'35.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.33%  "
$ws.Range("D35").Value = "'2.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").Value = "'49.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").Value = "'129.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").Value = "'16.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "'21.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'2.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").Value = "'2.060.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").Value = "'0.0317"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("E51").Value = "  +0.11%  "
